$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-26 06:36:16"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-26 06:36:12"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee29a162cb106bb16bd888368e03f645eb545fc0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bbc5d26e933ca2a2d4e4408c15a6bf4385a11794/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-26 06:36:16"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee29a162cb106bb16bd888368e03f645eb545fc0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bbc5d26e933ca2a2d4e4408c15a6bf4385a11794/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.17
